# "Some edits to design section." — thesis progress tracker workbook
#
# Sheet2 ("design" section rows + Score columns N:R):
#   - Add actual daily Score readings (N4:N7) that were missing, which in
#     turn ripple through the Actual Difference (P) / Goal Difference (Q)
#     calculations for the whole table (shared formulas recomputed).
#   - In the design checklist (rows 17-20) move a few items from the
#     "< 0.5" / "not started" columns into "Done".
#
# Also: the user had switched their active worksheet/selection to Sheet2
# (design section) before saving.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Score tracker: fill in actual recorded scores for a few days ---
$ws2.Range("N4").Value = 12.6
$ws2.Range("N5").Value = 13.3
$ws2.Range("N6").Value = 13.3
$ws2.Range("N7").Value = 15.3

# Q4 is its own (non-shared) formula; Q5:Q14 is a shared-formula block, so
# rewrite the whole block via R1C1 so the shared template updates for every
# row in one go (O[row] - N[row-1] instead of O[row] - O[row-1]).
$ws2.Range("Q4").Formula = "=O4-N3"
$ws2.Range("Q5:Q14").FormulaR1C1 = "=RC[-2]-R[-1]C[-3]"

# --- Design checklist: mark items Done that were previously in another
# progress bucket ---
$ws2.Range("D17").Value = 1
$ws2.Range("F17").Value = ""

$ws2.Range("D18").Value = 1
$ws2.Range("F18").Value = ""

$ws2.Range("D19").Value = 1
$ws2.Range("G19").Value = ""

$ws2.Range("D20").Value = 1
$ws2.Range("G20").Value = ""

# --- Active sheet / selection bookkeeping ---
# The author left off with the cursor on Sheet1!E27 before switching over
# to work in Sheet2, which is the tab that was active/selected on save.
[void]$ws1.Range("E27").Select()
[void]$ws2.Activate()
[void]$ws2.Range("G23").Select()
